$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Examples & Info")

$ws.Range("AB2").Value = 'Sequence of the I7 barcode (for 10x-Single Cell, you need enter four rows to supply the barcodes and add a "_a,_b,_c, _d" suffix to the Sample Type)'
$ws.Range("AB3").Value = "TTACCGAC"
$ws.Range("AB6").Value = "[ATGC]*"
